$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Subscript-3 unicode character used in D21 (0.0₃0944)
$sub3 = [string][char]0x2083

# Cells whose new values look numeric but must remain stored as text
$textCells = @("D5","D6","D7","D8","D9","D10","D11","D13","D17","D19","D20","D22","D23","D24","D26","D28","D29","D30","D31","D32","D33","D34","D37","D38","D39","D41","D42","D43","D46","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = '47.353.65'
$ws.Range("E2").Value = '  +3.16%  '
$ws.Range("D3").Value = '2.509.91'
$ws.Range("E3").Value = '  +2.58%  '
$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("D5").Value = '324.03'
$ws.Range("E5").Value = '  +1.10%  '
$ws.Range("D6").Value = '109.83'
$ws.Range("E6").Value = '  +5.70%  '
$ws.Range("D7").Value = '0.526'
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.539'
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("D10").Value = '39.28'
$ws.Range("E10").Value = '  +9.89%  '
$ws.Range("D11").Value = '0.0817'
$ws.Range("E11").Value = '  +1.68%  '
$ws.Range("D13").Value = '18.59'
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("E14").Value = '  +2.89%  '
$ws.Range("D15").Value = '2.903.13'
$ws.Range("E15").Value = '  +2.59%  '
$ws.Range("D16").Value = '2.511.68'
$ws.Range("E16").Value = '  +2.92%  '
$ws.Range("D17").Value = '0.862'
$ws.Range("E17").Value = '  +2.89%  '
$ws.Range("D18").Value = '47.316.99'
$ws.Range("E18").Value = '  +3.43%  '
$ws.Range("D19").Value = '12.92'
$ws.Range("E19").Value = '  +3.53%  '
$ws.Range("D20").Value = '6.70'
$ws.Range("E20").Value = '  +4.75%  '
$ws.Range("D21").Value = "0.0" + $sub3 + "0944"
$ws.Range("E21").Value = '  +1.36%  '
$ws.Range("D22").Value = '2.65'
$ws.Range("E22").Value = '  +13.15%  '
$ws.Range("D23").Value = '70.76'
$ws.Range("E23").Value = '  -0.66%  '
$ws.Range("D24").Value = '249.57'
$ws.Range("E24").Value = '  +1.52%  '
$ws.Range("E25").Value = '  +3.98%  '
$ws.Range("D26").Value = '26.15'
$ws.Range("E26").Value = '  +0.97%  '
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = '2.29'
$ws.Range("E28").Value = '  +4.77%  '
$ws.Range("D29").Value = '10.04'
$ws.Range("E29").Value = '  +3.64%  '
$ws.Range("D30").Value = '35.59'
$ws.Range("E30").Value = '  +5.75%  '
$ws.Range("D31").Value = '0.137'
$ws.Range("E31").Value = '  +6.89%  '
$ws.Range("D32").Value = '50.32'
$ws.Range("E32").Value = '  +2.26%  '
$ws.Range("D33").Value = '19.96'
$ws.Range("E33").Value = '  +0.53%  '
$ws.Range("D34").Value = '5.45'
$ws.Range("E34").Value = '  +2.07%  '
$ws.Range("E35").Value = '  +4.87%  '
$ws.Range("E36").Value = '  +0.25%  '
$ws.Range("D37").Value = '2.00'
$ws.Range("E37").Value = '  +6.05%  '
$ws.Range("D38").Value = '4.75'
$ws.Range("E38").Value = '  +4.95%  '
$ws.Range("D39").Value = '3.02'
$ws.Range("E39").Value = '  +3.44%  '
$ws.Range("E40").Value = '  +1.81%  '
$ws.Range("D41").Value = '122.45'
$ws.Range("E41").Value = '  -2.73%  '
$ws.Range("D42").Value = '2.24'
$ws.Range("E42").Value = '  -0.88%  '
$ws.Range("D43").Value = '21.23'
$ws.Range("E43").Value = '  +1.62%  '
$ws.Range("E44").Value = '  +2.45%  '
$ws.Range("D45").Value = '2.006.12'
$ws.Range("E45").Value = '  +2.60%  '
$ws.Range("D46").Value = '3.11'
$ws.Range("E46").Value = '  +5.45%  '
$ws.Range("E47").Value = '  -1.66%  '
$ws.Range("B48").Value = 'Stacks'
$ws.Range("C48").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D48").Value = '1.78'
$ws.Range("E48").Value = '  -3.43%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '9.09'
$ws.Range("E49").Value = '  +0.02%  '
$ws.Range("D50").Value = '5.26'
$ws.Range("E50").Value = '  +6.85%  '
$ws.Range("D51").Value = '78.41'
$ws.Range("E51").Value = '  +1.19%  '
